{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block that\n// followed the \"LOQ4240: ...\" requirements paragraph, along with the blank\n// paragraph that separated it from that line. The blank paragraph that used\n// to sit between the copyright line and the page-break paragraph is kept in\n// place (it now directly follows the \"LOQ4240...\" paragraph).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its text so the script does not depend on\n// a brittle, hard-coded paragraph index.\nconst anchorText = \"LOQ4240: Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II (Requisito fraco)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// The three paragraphs immediately after the anchor are the ones being\n// removed: a blank paragraph, the \"Ver no Jupiter...\" line, and the\n// \"\u00a9 2020 ...\" copyright line.\nconst expectedRemoved = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nfor (let k = 0; k < expectedRemoved.length; k++) {\n  const idx = anchorIndex + 1 + k;\n  if (idx >= paragraphs.items.length || paragraphs.items[idx].text !== expectedRemoved[k]) {\n    throw new Error(\"Unexpected document structure near anchor paragraph; aborting to avoid deleting the wrong content.\");\n  }\n}\n\n// Delete from the last to the first so earlier indices stay valid.\nfor (let k = expectedRemoved.length - 1; k >= 0; k--) {\n  paragraphs.items[anchorIndex + 1 + k].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block that\n# followed the \"LOQ4240: ...\" requirements paragraph, along with the blank\n# paragraph that separated it from that line. The blank paragraph that used\n# to sit between the copyright line and the page-break paragraph is kept in\n# place (it now directly follows the \"LOQ4240...\" paragraph).\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its text so the script does not depend on\n# a brittle, hard-coded paragraph index.\n$anchorText = \"LOQ4240: Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II (Requisito fraco)\"\n\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd(\"`r\")\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# The three paragraphs immediately after the anchor are the ones being\n# removed: a blank paragraph, the \"Ver no Jupiter...\" line, and the\n# \"\u00a9 2020 ...\" copyright line.\n$expectedRemoved = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\nfor ($k = 0; $k -lt $expectedRemoved.Count; $k++) {\n    $idx = $anchorIndex + 1 + $k\n    $t = $d.Paragraphs($idx).Range.Text.TrimEnd(\"`r\")\n    if ($t -ne $expectedRemoved[$k]) {\n        throw \"Unexpected document structure near anchor paragraph; aborting to avoid deleting the wrong content.\"\n    }\n}\n\n# Delete from the last to the first so earlier indices stay valid.\nfor ($k = $expectedRemoved.Count - 1; $k -ge 0; $k--) {\n    $idx = $anchorIndex + 1 + $k\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
